# Update countries & provincias Spain
# Applies the data refresh captured in the commit:
#  - Updated "Datos actualizados" timestamp in A1
#  - Peru overtakes Mexico in the ranking (new Peru numbers, Mexico keeps
#    its previous values but drops one row)
#  - Bolivia overtakes Egipto in the ranking (new Bolivia numbers, Egipto
#    keeps its previous values but drops one row)
#  - Montserrat overtakes Islas Malvinas in the ranking (Montserrat keeps
#    its previous values and moves up one row, Islas Malvinas keeps its
#    previous values and drops one row)
#  - Straightforward numeric refreshes for Kazajistan, Belgica, Venezuela,
#    Camboya, San Martin (Parte Holandesa) and Islas Turcas y Caicos

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-Row($r, $values) {
    $cols = @("A","B","C","D","E","F","G","H")
    for ($i = 0; $i -lt $values.Length; $i++) {
        $ws.Range($cols[$i] + $r).Value = $values[$i]
    }
}

# Updated timestamp
$ws.Range("A1").Value = "Datos actualizados a 14 de Agosto de 2020 a las 05:07"

# Row 9/10: Peru moves above Mexico
Set-Row 9  @("Peru", 507996, 0, 348006, 138277, 0, 0, 21713)
Set-Row 10 @("Mexico", 505751, 7371, 341507, 108951, 0, 627, 55293)

# Row 29: Kazajistan refreshed numbers (no reorder)
Set-Row 29 @("Kazajistan", 101848, 476, 76756, 23823, 0, 0, 1269)

# Row 31/32: Bolivia moves above Egipto
Set-Row 31 @("Bolivia", 96459, 1388, 33720, 58855, 0, 57, 3884)
Set-Row 32 @("Egipto", 96108, 0, 56890, 34111, 0, 0, 5107)

# Row 40: Belgica refreshed numbers (no reorder)
Set-Row 40 @("Belgica", 76191, 544, 17913, 48362, 0, 16, 9916)

# Row 63: Venezuela refreshed numbers (no reorder)
Set-Row 63 @("Venezuela", 30369, 0, 21385, 8725, 0, 0, 259)

# Row 180: Camboya refreshed numbers (no reorder)
Set-Row 180 @("Camboya", 273, 1, 225, 48, 0, 0, 0)

# Row 182: San Martin (Parte Holandesa) refreshed numbers (no reorder)
Set-Row 182 @("San Martin (Parte Holandesa)", 263, 15, 102, 144, 0, 0, 17)

# Row 183: Islas Turcas y Caicos refreshed numbers (no reorder)
Set-Row 183 @("Islas Turcas y Caicos", 258, 17, 52, 204, 0, 0, 2)

# Row 213/214: Montserrat moves above Islas Malvinas
Set-Row 213 @("Montserrat", 13, 0, 12, 0, 0, 0, 1)
Set-Row 214 @("Islas Malvinas", 13, 0, 13, 0, 0, 0, 0)
